$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks before rewriting cell values/content so that
# stale hyperlink associations (pointing at the old cells) do not linger.
$ws.Hyperlinks.Delete()

# --- Rebuild grid rows 2-28 (row 1 header is unchanged) ---
$ws.Range("A2").Value = "goto"
$ws.Range("B2").Value = "https://www.bluenile.com/jewelry/necklaces/lab-grown-diamond-cushion-cut-solitaire-pendant-in-14k-white-gold-1-2-ct-tw-f-g-vs2-si1-item-202314"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 9000

$ws.Range("A3").Value = "scroll"
$ws.Range("B3").Value = "Ships by"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 5000

$ws.Range("A4").Value = "click"
$ws.Range("B4").Value = "ADD TO CART button"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 2000

$ws.Range("A5").Value = "goto"
$ws.Range("B5").Value = "https://www.bluenile.com/shopping-cart"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 9000

$ws.Range("A6").Value = "wairfortext"
$ws.Range("B6").Value = "Summary"
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

$ws.Range("A7").Value = "scroll"
$ws.Range("B7").Value = "We Accept"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 1000
$ws.Range("E7").Value = 5000

$ws.Range("A8").Value = "click"
$ws.Range("B8").Value = "Checkout button"
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 1000
$ws.Range("E8").Value = 5000

$ws.Range("A9").Value = "wairfortext"
$ws.Range("B9").Value = "Please provide an email address"
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()

$ws.Range("A10").Value = "click"
$ws.Range("B10").Value = "Email Address input field"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 1000
$ws.Range("E10").Value = 2000

$ws.Range("A11").Value = "type"
$ws.Range("B11").Value = "Email Address input field"
$ws.Range("C11").Value = "mellina@gmail.com"
$ws.Range("D11").Value = 1000
$ws.Range("E11").Value = 2000

$ws.Range("A12").Value = "click"
$ws.Range("B12").Value = "Continue button"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = 1000
$ws.Range("E12").Value = 5000

$ws.Range("A13").Value = "wairfortext"
$ws.Range("B13").Value = "First Name"
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

$ws.Range("A14").Value = "click"
$ws.Range("B14").Value = "First Name input field"
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 2000

$ws.Range("A15").Value = "type"
$ws.Range("B15").Value = "First Name input field"
$ws.Range("C15").Value = "Mellina"
$ws.Range("D15").Value = 1000
$ws.Range("E15").Value = 2000

$ws.Range("A16").Value = "click"
$ws.Range("B16").Value = "Last Name input field"
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = 1000
$ws.Range("E16").Value = 2000

$ws.Range("A17").Value = "type"
$ws.Range("B17").Value = "Last Name input field"
$ws.Range("C17").Value = "James"
$ws.Range("D17").Value = 1000
$ws.Range("E17").Value = 2000

$ws.Range("A18").Value = "click"
$ws.Range("B18").Value = "Enter Address Manually button"
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = 1000
$ws.Range("E18").Value = 2000

$ws.Range("A19").Value = "scroll"
$ws.Range("B19").Value = "Use same address for billing"
$ws.Range("C19").ClearContents()
$ws.Range("D19").Value = 1000
$ws.Range("E19").Value = 3000

$ws.Range("A20").Value = "filldata"
$ws.Range("B20").Value = "input#address-finder"
$ws.Range("C20").Value = "3710 Pio Pico St"
$ws.Range("D20").Value = 1000
$ws.Range("E20").Value = 2000

$ws.Range("A21").Value = "presskey"
$ws.Range("B21").Value = "input#address-finder"
$ws.Range("C21").Value = "Enter"
$ws.Range("D21").Value = 1000
$ws.Range("E21").Value = 2000

$ws.Range("A22").Value = "scroll"
$ws.Range("B22").Value = "Continue button"
$ws.Range("C22").ClearContents()
$ws.Range("D22").Value = 1000
$ws.Range("E22").Value = 2000

$ws.Range("A23").Value = "click"
$ws.Range("B23").Value = "Phone Number input field"
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = 1000
$ws.Range("E23").Value = 2000

$ws.Range("A24").Value = "type"
$ws.Range("B24").Value = "Phone Number input field"
$ws.Range("C24").Value = 6142273098
$ws.Range("D24").Value = 1000
$ws.Range("E24").Value = 2000

$ws.Range("A25").Value = "clickloc"
$ws.Range("B25").Value = "button[data-qa=`"continue_button-shipping_and_billing_step-checkout_page`"]"
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = 1000
$ws.Range("E25").Value = 5000

$ws.Range("A26").Value = "scroll"
$ws.Range("B26").Value = "Contact information"
$ws.Range("C26").ClearContents()
$ws.Range("D26").Value = 1000
$ws.Range("E26").Value = 5000

$ws.Range("A27").Value = "assert"
$ws.Range("B27").Value = "#email"
$ws.Range("C27").Value = "mellina@gmail.com"
$ws.Range("D27").Value = 1000
$ws.Range("E27").Value = 2000

$ws.Range("A28").Value = "assert"
$ws.Range("B28").Value = "#phone"
$ws.Range("C28").Value = 6142273098
$ws.Range("D28").Value = 1000
$ws.Range("E28").Value = 2000


# --- Hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.bluenile.com/jewelry/necklaces/lab-grown-diamond-cushion-cut-solitaire-pendant-in-14k-white-gold-1-2-ct-tw-f-g-vs2-si1-item-202314") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:mellina@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C27"), "mailto:mellina@gmail.com") | Out-Null

# --- Sheet view updates (scroll position / selection) ---
$ws.Range("A24").Select()
$excel.ActiveWindow.ScrollRow = 7

$wb.Save()
